$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count()

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p -eq "System" -or $p -eq "system") {
                $hasSystem = $true
            }
        }
        if ($hasSystem) {
            $n = $parts.Count
            $rev = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $rev += $parts[$i]
            }
            $newVal = [string]::Join(", ", $rev)
            $cell.Value = $newVal
        }
    }
}
